# feat: add 2022-Q4 data
#
# Target layout:
#   Sheet "总计"    (unchanged position 1, gets a new row for 2021-Q1 and an
#                    updated row for 2022-Q4)
#   Sheet "2022-Q4" (new data, takes over the old "2021-Q1" sheet slot/position 2)
#   Sheet "2021-Q1" (new sheet, position 3, holds the data the old "2021-Q1"
#                    sheet used to have)

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)      # "总计"
$wsOld   = $wb.Worksheets.Item(2)      # currently "2021-Q1", holds the old fund table

# ------------------------------------------------------------------
# 1. Create the new sheet right after the existing one (temporary
#    name, since "2021-Q1" is still taken by $wsOld) and copy the
#    *formatting* of the old sheet's header/A-column into it before
#    we start overwriting the old sheet's own content.
# ------------------------------------------------------------------
$wsNew = $wb.Worksheets.Add($null, $wsOld)
$wsNew.Name = "2021-Q1-new"

$wsOld.Range("B1:H1").Copy()
$wsNew.Range("B1:H1").PasteSpecial(-4122)

$wsOld.Range("A2:A3").Copy()
$wsNew.Range("A2:A3").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2. Populate the new "2021-Q1" sheet with the data that used to live
#    on the old sheet (fund table for 华夏新起点灵活配置混合). The
#    numeric-looking values must stay plain text (no leading-zero /
#    trailing-zero loss), so force a text format before writing them
#    and drop back to the default "Normal" style afterwards so no
#    extra cell style gets created.
# ------------------------------------------------------------------
$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金金额"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

$wsNew.Range("A2").Value = 0
$wsNew.Range("C2").Value = "华夏新起点灵活配置混合A"
$wsNew.Range("H2").Value = 7

$wsNew.Range("A3").Value = 1
$wsNew.Range("C3").Value = "华夏新起点灵活配置混合C"
$wsNew.Range("H3").Value = 7

$wsNewText = $wsNew.Range("B2:G3")
$wsNewText.NumberFormat = "@"
$wsNew.Range("B2").Value = "002604"
$wsNew.Range("D2").Value = "0.63"
$wsNew.Range("E2").Value = "30.98"
$wsNew.Range("F2").Value = "2.03"
$wsNew.Range("G2").Value = "0.0128"
$wsNew.Range("B3").Value = "008213"
$wsNew.Range("D3").Value = "0.01"
$wsNew.Range("E3").Value = "30.98"
$wsNew.Range("F3").Value = "2.03"
$wsNew.Range("G3").Value = "0.0002"
$wsNewText.Style = "Normal"

# ------------------------------------------------------------------
# 3. Turn the old sheet into "2022-Q4": give it the header/A-column
#    formatting used on the "总计" sheet, then fill in the new data.
#    Rename it first so the temporary "2021-Q1-new" name can become
#    the real "2021-Q1".
# ------------------------------------------------------------------
$wsOld.Name = "2022-Q4"
$wsNew.Name = "2021-Q1"

$wsTotal.Range("B1:D1").Copy()
$wsOld.Range("B1:D1").PasteSpecial(-4122)
$wsTotal.Range("D1").Copy()
$wsOld.Range("E1:H1").PasteSpecial(-4122)

$wsTotal.Range("A2").Copy()
$wsOld.Range("A2:A3").PasteSpecial(-4122)

$wsOld.Range("B1").Value = "基金代码"
$wsOld.Range("C1").Value = "基金名称"
$wsOld.Range("D1").Value = "基金规模"
$wsOld.Range("E1").Value = "股票总仓位"
$wsOld.Range("F1").Value = "仓位占比"
$wsOld.Range("G1").Value = "持有市值(亿元)"
$wsOld.Range("H1").Value = "仓位排名"

$wsOld.Range("A2").Value = 0
$wsOld.Range("C2").Value = "招商量化精选股票A"
$wsOld.Range("H2").Value = 4

$wsOld.Range("A3").Value = 1
$wsOld.Range("C3").Value = "招商量化精选股票C"
$wsOld.Range("H3").Value = 4

$wsOldText = $wsOld.Range("B2:G3")
$wsOldText.NumberFormat = "@"
$wsOld.Range("B2").Value = "001917"
$wsOld.Range("D2").Value = "5.91"
$wsOld.Range("E2").Value = "94.08"
$wsOld.Range("F2").Value = "1.45"
$wsOld.Range("G2").Value = "0.0857"
$wsOld.Range("B3").Value = "007950"
$wsOld.Range("D3").Value = "5.28"
$wsOld.Range("E3").Value = "94.08"
$wsOld.Range("F3").Value = "1.45"
$wsOld.Range("G3").Value = "0.0766"
$wsOldText.Style = "Normal"

# ------------------------------------------------------------------
# 4. Update the "总计" summary sheet: row 2 now reports 2022-Q4, and a
#    new row 3 is added (copying row 2's A-column style) for 2021-Q1.
# ------------------------------------------------------------------
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.16

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q1"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.01
